# Report_Components.xlsx edit
# - Rename OPTIONAL_4's table label from "Annotation Table" to
#   "Refined Annotation Table".
# - Add a new OPTIONAL_5 row ("Artifact Annotation Table") right after it,
#   mirroring the existing OPTIONAL_4 row's layout/format.
# - Normalize the stray duplicate header style on A2:B2 to match the
#   identical style already used by the C1:E1 header cells.
# - Leave selection on the first empty row below the table (A18), matching
#   where Excel would land the cursor after appending a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2:B2 carried a near-duplicate of the C1:E1 header style (only differing
# by an inconsequential applyBorder flag). Re-pasting the C1 format over
# A2:B2 collapses them onto the same style entry.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null

# Seed row 17 with row 16's formatting before filling in its values so the
# new row matches the table's existing look (fonts, fills, alignment).
$ws.Range("A16:F16").Copy() | Out-Null
$ws.Range("A17:F17").PasteSpecial(-4122) | Out-Null

# Rename the existing "Annotation Table" entry.
$ws.Range("B16").Value = "Refined Annotation Table"

# Populate the new OPTIONAL_5 row.
$ws.Range("A17").Value = "OPTIONAL_5"
$ws.Range("B17").Value = "Artifact Annotation Table"
$ws.Range("C17").Value = "y"
$ws.Range("D17").Value = "y"
$ws.Range("E17").Value = "y"
$ws.Range("F17").Value = "optional"

# Match the post-edit active selection recorded in the saved workbook.
$ws.Range("A18").Select() | Out-Null
